$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Fix text typos (stray embedded spaces / full-width commas) in existing values
$ws.Cells.Item(2,2).Value = "太平洋電線電纜股份有限公司"
$ws.Cells.Item(10,2).Value = "太平洋電線電纜股份有限公司"
$ws.Cells.Item(7,2).Value = "春雨開發股份有限公司（原正華）"
$ws.Cells.Item(12,2).Value = "台灣土地開發股份有限公司"
$ws.Cells.Item(11,4).Value = "2758"
$ws.Cells.Item(8,7).Value = "4870"

# Insert new "property_category" column before the "date" column and fill it with "stock"
$ws.Columns.Item(8).Insert()
$ws.Cells.Item(1,8).Value = "property_category"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r,8).Value = "stock"
}
